# Add files via upload
# Rebuild the "100_2" confirmations sheet: break each summary bucket
# (New nominations / Carryover nominations / Confirmed / Unconfirmed /
# Withdrawn) out by branch name, drop the stand-alone "Summary" header
# row, and reorder/relabel the final four grand-total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the "Summary" header row (old row 38). Excel shifts every
#    row below it up by one, so the four grand totals that used to sit
#    at rows 39-43 land on rows 38-42 - giving the final 42-row sheet.
# ---------------------------------------------------------------------
$ws.Rows(38).Delete()

# ---------------------------------------------------------------------
# 2. Relabel column A so each bucket row spells out which branch it
#    belongs to, e.g. "     New nominations" -> "     Civilian, New
#    nominations". Header/category rows (Civilian, Army, Navy, ...)
#    just gain a trailing space. Rows that don't change text are left
#    alone.
# ---------------------------------------------------------------------
$labels = @{
    6  = "Civilian "
    7  = "     Civilian, New nominations"
    8  = "     Civilian, Carryover nominations"
    9  = "     Civilian, Confirmed  "
    10 = "     Civilian, Unconfirmed "
    11 = "     Civilian, Withdrawn "
    13 = "     Civilian (lists), New nominations"
    14 = "     Civilian (lists), Carryover nominations"
    15 = "     Civilian (lists), Confirmed  "
    16 = "     Civilian (lists), Unconfirmed  "
    17 = "Air Force "
    18 = "     Air Force, New nominations"
    19 = "     Air Force, Carryover nominations"
    20 = "     Air Force, Confirmed "
    21 = "Army "
    22 = "     Army, New nominations"
    23 = "     Army, Carryover nominations"
    24 = "     Army, Confirmed "
    25 = "     Army, Unconfirmed"
    26 = "     Army, Withdrawn  "
    27 = "Navy "
    28 = "     Navy, New nominations"
    29 = "     Navy, Carryover nominations"
    30 = "     Navy, Confirmed "
    31 = "     Navy, Unconfirmed"
    32 = "     Navy, Withdrawn  "
    33 = "Marine Corps "
    34 = "     Marine Corps, New nominations"
    35 = "     Marine Corps, Carryover nominations"
    36 = "     Marine Corps, Confirmed"
    37 = "     Marine Corps, Unconfirmed"
    38 = "Total new nominations"
    39 = "Total carryover nominations"
    40 = "Total confirmed "
}

foreach ($r in $labels.Keys) {
    $ws.Cells.Item($r, 1).Value = $labels[$r]
}

# ---------------------------------------------------------------------
# 3. The last two grand-total rows (39 "Total new nominations" and 38
#    "Total carryover nominations" - now 38/39 after the delete above)
#    swap which number belongs to which label, so rewrite those two
#    values explicitly along with the other figures that moved/changed.
# ---------------------------------------------------------------------
$ws.Cells.Item(38, 2).Value = 37264
$ws.Cells.Item(39, 2).Value = 5494

# ---------------------------------------------------------------------
# 4. Column-wide alignment: labels left-aligned, values right-aligned.
# ---------------------------------------------------------------------
$ws.Columns("A").HorizontalAlignment = -4131
$ws.Columns("B").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 5. Subtotal rows (each branch's "Confirmed" line, plus the three
#    grand totals at the bottom) use the thousands-separator number
#    format instead of plain General.
# ---------------------------------------------------------------------
$subtotalRows = @(15, 20, 24, 30, 36, 38, 39, 40)
foreach ($r in $subtotalRows) {
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"
}

# ---------------------------------------------------------------------
# 6. Army/Navy/Marine Corps "Withdrawn"/"Unconfirmed" rows that no
#    longer carry a figure - remove the value cell entirely (not just
#    its content) so column B is blank for that row, same as row 27/32
#    which never had a number to begin with. Done last so the
#    column-wide formatting above doesn't resurrect an empty styled
#    cell here.
# ---------------------------------------------------------------------
$ws.Cells.Item(26, 2).Clear()
$ws.Cells.Item(27, 2).Clear()
$ws.Cells.Item(37, 2).Clear()

# ---------------------------------------------------------------------
# 7. Cosmetic: match the saved selection from the authored workbook.
# ---------------------------------------------------------------------
$ws.Range("A16").Select()
